$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:H headers to E:I
$ws.Range("D1").EntireColumn.Insert()

# Set the new column D header text
$ws.Range("D1").Value = "Status"

# Copy the header style/format from the neighboring C1 cell onto the new D1 cell
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# The month label that used to be "Nov_2025" shifted from F1 to G1; update it to "Oct_2025"
$ws.Range("G1").Value = "Oct_2025"
